$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add two new rows of recursion practice data (row 21 and row 22)
$ws.Range("A21").Value = 45838
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 13

$ws.Range("A22").Value = 45839
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 14

# Update the active selection to match the author's final cursor position
$ws.Range("C23").Select()
